# Nightly cryptos-list refresh (GitHub Actions bot).
# Updates the "Price" (D) and "Volume(1h)" (E) columns on Sheet1 for the
# coin rows that moved since the last run. Values are written as literal
# text (leading apostrophe + style reset) so numeric-looking prices such
# as "544.73" or thousand-dotted "58.359.12" stay text, matching the sheet's
# existing inline-string convention instead of being auto-coerced to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'58.359.12"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.96%  "
$ws.Range("D3").Value = "'2.280.32"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -5.65%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'544.73"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.47%  "
$ws.Range("D6").Value = "'130.60"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.73%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  -2.96%  "
$ws.Range("D9").Value = "'2.279.10"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -5.65%  "
$ws.Range("E10").Value = "  -3.84%  "
$ws.Range("E11").Value = "  -2.87%  "
$ws.Range("E12").Value = "  +0.03%  "
$ws.Range("E13").Value = "  -5.33%  "
$ws.Range("D14").Value = "'23.56"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -5.32%  "
$ws.Range("D15").Value = "'2.686.10"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -5.76%  "
$ws.Range("D16").Value = "'58.338.75"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.83%  "
$ws.Range("E17").Value = "  -3.46%  "
$ws.Range("D18").Value = "'2.282.87"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -5.84%  "
$ws.Range("E19").Value = "  -6.02%  "
$ws.Range("E20").Value = "  -4.28%  "
$ws.Range("D21").Value = "'313.75"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.36%  "
$ws.Range("E22").Value = "  -4.60%  "
$ws.Range("E23").Value = "  +0.03%  "
$ws.Range("D24").Value = "'62.86"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.78%  "
$ws.Range("E25").Value = "  -3.97%  "
$ws.Range("E26").Value = "  -0.07%  "
$ws.Range("D27").Value = "'8.09"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -6.90%  "
$ws.Range("E28").Value = "  -6.35%  "
$ws.Range("E29").Value = "  -1.17%  "
$ws.Range("D30").Value = "'170.50"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.30%  "
$ws.Range("D31").Value = "'0.0₃0721"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -6.49%  "
$ws.Range("E32").Value = "  -0.08%  "
$ws.Range("E33").Value = "  -5.65%  "
$ws.Range("D34").Value = "'0.383"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.90%  "
$ws.Range("E36").Value = "  -4.00%  "
$ws.Range("E37").Value = "  +0.08%  "
$ws.Range("E38").Value = "  -5.12%  "
$ws.Range("D39").Value = "'3.93"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -6.05%  "
$ws.Range("D40").Value = "'37.88"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.34%  "
$ws.Range("E41").Value = "  -5.47%  "
$ws.Range("D42").Value = "'294.68"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -9.97%  "
$ws.Range("D43").Value = "'140.58"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.33%  "
$ws.Range("E44").Value = "  -5.70%  "
$ws.Range("D45").Value = "'0.0946"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.83%  "
$ws.Range("E46").Value = "  -3.45%  "
$ws.Range("E47").Value = "  -3.77%  "
$ws.Range("D48").Value = "'18.31"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -8.51%  "
$ws.Range("E49").Value = "  -4.12%  "
$ws.Range("E50").Value = "  -5.33%  "
$ws.Range("D51").Value = "'11.01"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.31%  "
